# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row = new value }
$updates = @{
    "展览" = @{
        2  = 15056
        3  = 19184
        5  = 145
        13 = 60
        22 = 8010
        29 = 6076
        32 = 174
        35 = 5476
        36 = 894
    }
    "全部类型" = @{
        2  = 15056
        3  = 19184
        5  = 145
        13 = 60
        23 = 8010
        32 = 6076
        35 = 174
        38 = 5476
        39 = 894
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowValues = $updates[$sheetName]
    foreach ($row in $rowValues.Keys) {
        $ws.Range("F$row").Value = $rowValues[$row]
    }
}
